$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L: ratio of K (2*v0) to I (Vab-Vauf), added for rows 2, 5, 11 and 17
# (the "tragende" measurement rows of the Millikan evaluation table).
$ws.Range("L2").Formula  = "=K2/I2"
$ws.Range("L5").Formula  = "=K5/I5"
$ws.Range("L11").Formula = "=K11/I11"
$ws.Range("L17").Formula = "=K17/I17"

# Leave the selection where the author ended up after entering the last formula.
$ws.Range("L18").Select()
